$wb = $excel.ActiveWorkbook

# Remove the "vol min", "vol max", "Expected part" and "Left/right" (laterality /
# number-of-parts) columns from every check-protocol sheet that has them, the
# same way a user would in Excel: select the whole columns, then delete them.
$sheetNames = @("Clinical Structures", "opt structures", "couch_structures")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $ws.Range("C1:F1").EntireColumn.Select()
    $ws.Range("C1:F1").EntireColumn.Delete()
}
